$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 45219
$ws.Range("A15").NumberFormat = "d-mmm"
$ws.Range("B15").Value = "Internship"
$ws.Range("C15").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

$ws.Range("C16").Select()
